# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") gets its
#    table style switched from the custom "Table_0" style
#    ({1F1D5EE5-41B1-48B1-8AB7-D32B9299C325}) to the built-in style
#    {3C78B833-4933-4BBF-A746-D3784E042469}.
#
# 2) The deck's applied colour theme (the one wired to SlideMaster1, i.e.
#    the one that actually paints every slide) is switched from the
#    colourful "Integral / Red Violet" scheme to the plain default
#    "Office" scheme - this is what a user does from
#    Design tab -> Variants/Colors -> (pick a different theme colour set).

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{3C78B833-4933-4BBF-A746-D3784E042469}")

# --- 2) Theme colours -------------------------------------------------
# Index order exposed by ThemeColorScheme: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink. Values below are the stock
# PowerPoint "Office" theme colours (RGB() style integers).
$tcs = $slide.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
